$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pre_phyX (column F) values for rows whose source image files
# had their padding removed, per the fix described in the commit message.
$updates = @{
    2 = 0.01065292096219931
    3 = 0.008591065292096219
    5 = 0.008013937282229966
    6 = 0.01132075471698113
    7 = 0.01047297297297297
    8 = 0.009326424870466322
    9 = 0.009605488850771868
    10 = 0.008362369337979094
    12 = 0.009411764705882352
    14 = 0.00798611111111111
    15 = 0.01105354058721934
    16 = 0.01030927835051546
    17 = 0.01063464837049743
    18 = 0.01006944444444445
    19 = 0.008013937282229966
    20 = 0.007342657342657343
    22 = 0.007355516637478108
    23 = 0.01045296167247387
    24 = 0.008141592920353982
    26 = 0.01090289608177172
    27 = 0.009326424870466322
    28 = 0.007815275310834813
    29 = 0.00909090909090909
    30 = 0.01035058430717863
    33 = 0.006137184115523466
    34 = 0.004128440366972477
    35 = 0.006206896551724138
    36 = 0.006930693069306931
    37 = 0.00684931506849315
    38 = 0.00472972972972973
    39 = 0.00466786355475763
    42 = 0.009482758620689655
    44 = 0.01104972375690608
    47 = 0.009557522123893806
    48 = 0.007355516637478108
    49 = 0.01217391304347826
    51 = 0.01039861351819757
    52 = 0.008235294117647059
    53 = 0.007381370826010545
    54 = 0.01218274111675127
    55 = 0.008981001727115715
    56 = 0.008445945945945945
    58 = 0.008156028368794326
    59 = 0.00684931506849315
    61 = 0.01030195381882771
    62 = 0.008726003490401396
    63 = 0.005385996409335727
    64 = 0.009605488850771868
    65 = 0.009491525423728813
    66 = 0.007521367521367522
    67 = 0.01172413793103448
    69 = 0.008510638297872341
    70 = 0.008981001727115715
    71 = 0.009863945578231292
    72 = 0.008156028368794326
    73 = 0.0101010101010101
    74 = 0.007719298245614034
    76 = 0.0101010101010101
    77 = 0.009294320137693631
    78 = 0.01206030150753769
    79 = 0.007725321888412016
    80 = 0.01154499151103565
    81 = 0.01105354058721934
    82 = 0.01120840630472855
    84 = 0.009948542024013723
    85 = 0.009215017064846417
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}

Write-Output "Updated $($updates.Count) cells in column F"